# Apply targeted numeric corrections to the robustness-table cells.
# Each entry is (row, column) in the single document table (1-indexed,
# as used by Word's Table.Cell(row, col)). We verify the existing value
# before overwriting it so the script fails loudly if the target drifts.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$c = $t.Cell(3, 5)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "0.133") {
    Write-Output "MISMATCH at row 3 col 5: got [$cur] expected [0.133]"
} else {
    $c.Range.Text = "0.134"
}

$c = $t.Cell(4, 4)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "[0.075,0.189]") {
    Write-Output "MISMATCH at row 4 col 4: got [$cur] expected [[0.075,0.189]]"
} else {
    $c.Range.Text = "[0.085,0.178]"
}

$c = $t.Cell(4, 5)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "[0.058,0.208]") {
    Write-Output "MISMATCH at row 4 col 5: got [$cur] expected [[0.058,0.208]]"
} else {
    $c.Range.Text = "[0.082,0.187]"
}

$c = $t.Cell(7, 4)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "0.099") {
    Write-Output "MISMATCH at row 7 col 4: got [$cur] expected [0.099]"
} else {
    $c.Range.Text = "0.090"
}

$c = $t.Cell(7, 5)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "0.110") {
    Write-Output "MISMATCH at row 7 col 5: got [$cur] expected [0.110]"
} else {
    $c.Range.Text = "0.082"
}

$c = $t.Cell(8, 4)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "[0.006,0.193]") {
    Write-Output "MISMATCH at row 8 col 4: got [$cur] expected [[0.006,0.193]]"
} else {
    $c.Range.Text = "[0.033,0.147]"
}

$c = $t.Cell(8, 5)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "[0.015,0.204]") {
    Write-Output "MISMATCH at row 8 col 5: got [$cur] expected [[0.015,0.204]]"
} else {
    $c.Range.Text = "[0.024,0.141]"
}

$c = $t.Cell(11, 5)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "0.089") {
    Write-Output "MISMATCH at row 11 col 5: got [$cur] expected [0.089]"
} else {
    $c.Range.Text = "0.084"
}

$c = $t.Cell(12, 5)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "[0.021,0.157]") {
    Write-Output "MISMATCH at row 12 col 5: got [$cur] expected [[0.021,0.157]]"
} else {
    $c.Range.Text = "[0.018,0.150]"
}

$c = $t.Cell(15, 4)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "0.254") {
    Write-Output "MISMATCH at row 15 col 4: got [$cur] expected [0.254]"
} else {
    $c.Range.Text = "0.252"
}

$c = $t.Cell(16, 4)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "[0.135,0.374]") {
    Write-Output "MISMATCH at row 16 col 4: got [$cur] expected [[0.135,0.374]]"
} else {
    $c.Range.Text = "[0.132,0.373]"
}

$c = $t.Cell(19, 4)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "0.204") {
    Write-Output "MISMATCH at row 19 col 4: got [$cur] expected [0.204]"
} else {
    $c.Range.Text = "0.100"
}

$c = $t.Cell(19, 5)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "0.133") {
    Write-Output "MISMATCH at row 19 col 5: got [$cur] expected [0.133]"
} else {
    $c.Range.Text = "0.112"
}

$c = $t.Cell(20, 4)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "[-0.111,0.520]") {
    Write-Output "MISMATCH at row 20 col 4: got [$cur] expected [[-0.111,0.520]]"
} else {
    $c.Range.Text = "[0.056,0.144]"
}

$c = $t.Cell(20, 5)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "[0.061,0.206]") {
    Write-Output "MISMATCH at row 20 col 5: got [$cur] expected [[0.061,0.206]]"
} else {
    $c.Range.Text = "[0.034,0.190]"
}

$c = $t.Cell(23, 4)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "0.064") {
    Write-Output "MISMATCH at row 23 col 4: got [$cur] expected [0.064]"
} else {
    $c.Range.Text = "0.091"
}

$c = $t.Cell(23, 5)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "0.131") {
    Write-Output "MISMATCH at row 23 col 5: got [$cur] expected [0.131]"
} else {
    $c.Range.Text = "0.134"
}

$c = $t.Cell(24, 4)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "[0.006,0.122]") {
    Write-Output "MISMATCH at row 24 col 4: got [$cur] expected [[0.006,0.122]]"
} else {
    $c.Range.Text = "[0.043,0.139]"
}

$c = $t.Cell(24, 5)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "[0.057,0.206]") {
    Write-Output "MISMATCH at row 24 col 5: got [$cur] expected [[0.057,0.206]]"
} else {
    $c.Range.Text = "[0.081,0.186]"
}

$c = $t.Cell(27, 4)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "0.201") {
    Write-Output "MISMATCH at row 27 col 4: got [$cur] expected [0.201]"
} else {
    $c.Range.Text = "0.172"
}

$c = $t.Cell(28, 4)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "[0.137,0.264]") {
    Write-Output "MISMATCH at row 28 col 4: got [$cur] expected [[0.137,0.264]]"
} else {
    $c.Range.Text = "[0.122,0.222]"
}

$c = $t.Cell(28, 5)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "[0.059,0.211]") {
    Write-Output "MISMATCH at row 28 col 5: got [$cur] expected [[0.059,0.211]]"
} else {
    $c.Range.Text = "[0.082,0.189]"
}

$c = $t.Cell(31, 5)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "0.132") {
    Write-Output "MISMATCH at row 31 col 5: got [$cur] expected [0.132]"
} else {
    $c.Range.Text = "0.133"
}

$c = $t.Cell(32, 4)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "[0.074,0.188]") {
    Write-Output "MISMATCH at row 32 col 4: got [$cur] expected [[0.074,0.188]]"
} else {
    $c.Range.Text = "[0.084,0.177]"
}

$c = $t.Cell(32, 5)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "[0.057,0.207]") {
    Write-Output "MISMATCH at row 32 col 5: got [$cur] expected [[0.057,0.207]]"
} else {
    $c.Range.Text = "[0.080,0.186]"
}

$c = $t.Cell(35, 4)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "0.135") {
    Write-Output "MISMATCH at row 35 col 4: got [$cur] expected [0.135]"
} else {
    $c.Range.Text = "0.134"
}

$c = $t.Cell(35, 5)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "0.137") {
    Write-Output "MISMATCH at row 35 col 5: got [$cur] expected [0.137]"
} else {
    $c.Range.Text = "0.138"
}

$c = $t.Cell(36, 4)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "[0.078,0.193]") {
    Write-Output "MISMATCH at row 36 col 4: got [$cur] expected [[0.078,0.193]]"
} else {
    $c.Range.Text = "[0.087,0.181]"
}

$c = $t.Cell(36, 5)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "[0.062,0.212]") {
    Write-Output "MISMATCH at row 36 col 5: got [$cur] expected [[0.062,0.212]]"
} else {
    $c.Range.Text = "[0.085,0.191]"
}

$c = $t.Cell(39, 4)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "0.163") {
    Write-Output "MISMATCH at row 39 col 4: got [$cur] expected [0.163]"
} else {
    $c.Range.Text = "0.158"
}

$c = $t.Cell(39, 5)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "0.172") {
    Write-Output "MISMATCH at row 39 col 5: got [$cur] expected [0.172]"
} else {
    $c.Range.Text = "0.167"
}

$c = $t.Cell(40, 4)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "[0.104,0.221]") {
    Write-Output "MISMATCH at row 40 col 4: got [$cur] expected [[0.104,0.221]]"
} else {
    $c.Range.Text = "[0.110,0.205]"
}

$c = $t.Cell(40, 5)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "[0.096,0.249]") {
    Write-Output "MISMATCH at row 40 col 5: got [$cur] expected [[0.096,0.249]]"
} else {
    $c.Range.Text = "[0.113,0.221]"
}

$c = $t.Cell(43, 4)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "0.123") {
    Write-Output "MISMATCH at row 43 col 4: got [$cur] expected [0.123]"
} else {
    $c.Range.Text = "0.122"
}

$c = $t.Cell(43, 5)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "0.124") {
    Write-Output "MISMATCH at row 43 col 5: got [$cur] expected [0.124]"
} else {
    $c.Range.Text = "0.122"
}

$c = $t.Cell(44, 4)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "[0.066,0.181]") {
    Write-Output "MISMATCH at row 44 col 4: got [$cur] expected [[0.066,0.181]]"
} else {
    $c.Range.Text = "[0.074,0.169]"
}

$c = $t.Cell(44, 5)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "[0.049,0.200]") {
    Write-Output "MISMATCH at row 44 col 5: got [$cur] expected [[0.049,0.200]]"
} else {
    $c.Range.Text = "[0.068,0.175]"
}

$c = $t.Cell(47, 4)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "0.117") {
    Write-Output "MISMATCH at row 47 col 4: got [$cur] expected [0.117]"
} else {
    $c.Range.Text = "0.120"
}

$c = $t.Cell(47, 5)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "0.109") {
    Write-Output "MISMATCH at row 47 col 5: got [$cur] expected [0.109]"
} else {
    $c.Range.Text = "0.120"
}

$c = $t.Cell(48, 4)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "[0.067,0.166]") {
    Write-Output "MISMATCH at row 48 col 4: got [$cur] expected [[0.067,0.166]]"
} else {
    $c.Range.Text = "[0.077,0.164]"
}

$c = $t.Cell(48, 5)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "[0.049,0.170]") {
    Write-Output "MISMATCH at row 48 col 5: got [$cur] expected [[0.049,0.170]]"
} else {
    $c.Range.Text = "[0.072,0.168]"
}

$c = $t.Cell(51, 4)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "0.109") {
    Write-Output "MISMATCH at row 51 col 4: got [$cur] expected [0.109]"
} else {
    $c.Range.Text = "0.112"
}

$c = $t.Cell(51, 5)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "0.102") {
    Write-Output "MISMATCH at row 51 col 5: got [$cur] expected [0.102]"
} else {
    $c.Range.Text = "0.109"
}

$c = $t.Cell(52, 4)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "[0.059,0.159]") {
    Write-Output "MISMATCH at row 52 col 4: got [$cur] expected [[0.059,0.159]]"
} else {
    $c.Range.Text = "[0.068,0.156]"
}

$c = $t.Cell(52, 5)
$cur = $c.Range.Text.TrimEnd([char]7, [char]13, [char]10)
if ($cur -ne "[0.041,0.162]") {
    Write-Output "MISMATCH at row 52 col 5: got [$cur] expected [[0.041,0.162]]"
} else {
    $c.Range.Text = "[0.060,0.158]"
}

